$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.118.62"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3
$ws.Range("D3").Value = "3.120.21"
$ws.Range("E3").Value = "  +0.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.05"
$ws.Range("E5").Value = "  -0.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.88"
$ws.Range("E6").Value = "  +1.13%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("E8").Value = "  -0.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.155"
$ws.Range("E9").Value = "  -0.24%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.42"
$ws.Range("E10").Value = "  -0.36%  "

# Row 11
$ws.Range("E11").Value = "  -0.70%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000249"
$ws.Range("E12").Value = "  -0.25%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.17"
$ws.Range("E13").Value = "  -0.59%  "

# Row 14
$ws.Range("E14").Value = "  -1.71%  "

# Row 15
$ws.Range("D15").Value = "3.636.08"
$ws.Range("E15").Value = "  +0.66%  "

# Row 16
$ws.Range("D16").Value = "67.111.43"
$ws.Range("E16").Value = "  +0.35%  "

# Row 18
$ws.Range("D18").Value = "3.119.24"
$ws.Range("E18").Value = "  +0.63%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.54"
$ws.Range("E19").Value = "  +1.92%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "489.75"
$ws.Range("E20").Value = "  +1.77%  "

# Row 21
$ws.Range("E21").Value = "  +5.82%  "

# Row 22
$ws.Range("E22").Value = "  -1.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.22"
$ws.Range("E23").Value = "  +0.28%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.22"
$ws.Range("E24").Value = "  +0.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  -3.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.38"
$ws.Range("E26").Value = "  +3.09%  "

# Row 27
$ws.Range("E27").Value = "  -0.07%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.91"
$ws.Range("E28").Value = "  -0.89%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.35"
$ws.Range("E29").Value = "  -2.10%  "

# Row 30
$ws.Range("E30").Value = "  -0.44%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.57"
$ws.Range("E31").Value = "  -0.80%  "

# Row 32
$ws.Range("E32").Value = "  -0.79%  "

# Row 33
$ws.Range("D33").Value = "0.0₃0949"
$ws.Range("E33").Value = "  -5.58%  "

# Row 34
$ws.Range("E34").Value = "  +0.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.87"
$ws.Range("E35").Value = "  -0.50%  "

# Row 36
$ws.Range("E36").Value = "  -1.70%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.16"
$ws.Range("E37").Value = "  -2.01%  "

# Row 38
$ws.Range("E38").Value = "  -3.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.310"
$ws.Range("E39").Value = "  -1.94%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.123"
$ws.Range("E40").Value = "  +1.41%  "

# Row 41
$ws.Range("E41").Value = "  -1.68%  "

# Row 42
$ws.Range("D42").Value = "2.819.34"
$ws.Range("E42").Value = "  -0.77%  "

# Row 43: dogwifhat -> Bittensor (with updated Price/Volume)
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "383.50"
$ws.Range("E43").Value = "  -0.58%  "

# Row 44: Bittensor -> dogwifhat (with updated Price/Volume)
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("E44").Value = "  -7.56%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0353"
$ws.Range("E45").Value = "  -2.38%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.60"
$ws.Range("E46").Value = "  +0.76%  "

# Row 47
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.87"
$ws.Range("E48").Value = "  -0.52%  "

# Row 49
$ws.Range("E49").Value = "  -1.24%  "

# Row 50
$ws.Range("E50").Value = "  -0.75%  "

# Row 51
$ws.Range("E51").Value = "  -1.91%  "

